$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = 484300
$ws.Range("E8").Value = 446000
$ws.Range("F8").Value = 386000
$ws.Range("G8").Value = 288600

$ws.Range("D9").Value = 331100
$ws.Range("E9").Value = 304800
$ws.Range("F9").Value = 259100
$ws.Range("G9").Value = 193200

$ws.Range("D10").Value = 153200
$ws.Range("E10").Value = 141200
$ws.Range("F10").Value = 127000
$ws.Range("G10").Value = 95500

$ws.Range("D15").Value = 8200
$ws.Range("E15").Value = 8000
$ws.Range("F15").Value = 5200
$ws.Range("G15").Value = 4800

$ws.Range("D17").Value = 501900
$ws.Range("E17").Value = 464600
$ws.Range("F17").Value = 402500
$ws.Range("G17").Value = 314900

$ws.Range("D18").Value = -17600
$ws.Range("E18").Value = -18600
$ws.Range("F18").Value = -16400
$ws.Range("G18").Value = -26200

$ws.Range("D20").Value = -9800
$ws.Range("E20").Value = -3100
$ws.Range("F20").Value = 5900

$ws.Range("D21").Value = -19400
$ws.Range("E21").Value = -14000
$ws.Range("F21").Value = -5300
$ws.Range("G21").Value = -20600

$ws.Range("D22").Value = 16300
$ws.Range("E22").Value = 17200
$ws.Range("F22").Value = 15000
$ws.Range("G22").Value = 12100

$ws.Range("D23").Value = -43700
$ws.Range("E23").Value = -38900
$ws.Range("F23").Value = -25500
$ws.Range("G23").Value = -37000

$ws.Range("D26").Value = -43700
$ws.Range("E26").Value = -38900
$ws.Range("F26").Value = -25500
$ws.Range("G26").Value = -37000

$ws.Range("D27").Value = -43500
$ws.Range("E27").Value = -38700
$ws.Range("F27").Value = -25300
$ws.Range("G27").Value = -36900

$ws.Range("D32").Value = 9800
$ws.Range("E32").Value = 3100
$ws.Range("F32").Value = -5900

$ws.Range("D33").Value = -43500
$ws.Range("E33").Value = -38700
$ws.Range("F33").Value = -25300
$ws.Range("G33").Value = -36900

$ws.Range("D35").Value = -43500
$ws.Range("E35").Value = -38700
$ws.Range("F35").Value = -25300
$ws.Range("G35").Value = -36900

$ws.Range("D41").Value = 101500
$ws.Range("E41").Value = 28500
$ws.Range("F41").Value = 63900
$ws.Range("G41").Value = 62100

$ws.Range("D43").Value = 49500
$ws.Range("E43").Value = 71900
$ws.Range("F43").Value = 86600
$ws.Range("G43").Value = 78200

$ws.Range("D44").Value = 117100
$ws.Range("E44").Value = 90200
$ws.Range("F44").Value = 75600
$ws.Range("G44").Value = 45600

$ws.Range("D45").Value = 17400
$ws.Range("E45").Value = 20800
$ws.Range("F45").Value = 14500
$ws.Range("G45").Value = 4600

$ws.Range("D46").Value = 285500
$ws.Range("E46").Value = 211400
$ws.Range("F46").Value = 240600
$ws.Range("G46").Value = 190600

$ws.Range("D47").Value = 18100
$ws.Range("E47").Value = 8500

$ws.Range("D48").Value = 18700
$ws.Range("E48").Value = 19000
$ws.Range("F48").Value = 17400
$ws.Range("G48").Value = 14000

$ws.Range("D49").Value = 29700
$ws.Range("E49").Value = 22500
$ws.Range("F49").Value = 13500
$ws.Range("G49").Value = 10400

$ws.Range("D52").Value = 31800
$ws.Range("E52").Value = 24100
$ws.Range("F52").Value = 14100
$ws.Range("G52").Value = 6000

$ws.Range("D54").Value = 383800
$ws.Range("E54").Value = 285500
$ws.Range("F54").Value = 285500
$ws.Range("G54").Value = 221000

$ws.Range("D57").Value = 93800
$ws.Range("E57").Value = 86000
$ws.Range("F57").Value = 66000
$ws.Range("G57").Value = 37300

$ws.Range("D58").Value = 27300
$ws.Range("E58").Value = 19500
$ws.Range("F58").Value = 26300
$ws.Range("G58").Value = 33100

$ws.Range("D59").Value = 83000
$ws.Range("E59").Value = 52600
$ws.Range("F59").Value = 41900
$ws.Range("G59").Value = 26600

$ws.Range("D60").Value = 204200
$ws.Range("E60").Value = 158100
$ws.Range("F60").Value = 134200
$ws.Range("G60").Value = 97000

$ws.Range("D61").Value = 46000
$ws.Range("E61").Value = 79800
$ws.Range("F61").Value = 59300
$ws.Range("G61").Value = 52900

$ws.Range("D62").Value = 9800
$ws.Range("E62").Value = 15800
$ws.Range("F62").Value = 17900
$ws.Range("G62").Value = 8300

$ws.Range("D66").Value = 259900
$ws.Range("E66").Value = 253800
$ws.Range("F66").Value = 211600
$ws.Range("G66").Value = 158600

$ws.Range("D72").Value = -217200
$ws.Range("E72").Value = -173700
$ws.Range("F72").Value = -134900
$ws.Range("G72").Value = -109600

$ws.Range("D76").Value = 123900
$ws.Range("E76").Value = 31700
$ws.Range("F76").Value = 73900
$ws.Range("G76").Value = 62400

$ws.Range("D81").Value = -43500
$ws.Range("E81").Value = -38700
$ws.Range("F81").Value = -25300
$ws.Range("G81").Value = -36900

$ws.Range("D83").Value = 8200
$ws.Range("E83").Value = 8000
$ws.Range("F83").Value = 5400
$ws.Range("G83").Value = 4400

$ws.Range("D89").Value = 6400
$ws.Range("E89").Value = -5300
$ws.Range("F89").Value = -5900
$ws.Range("G89").Value = -7000

$ws.Range("D91").Value = -2000
$ws.Range("E91").Value = -6500
$ws.Range("F91").Value = -5600
$ws.Range("G91").Value = -1400

$ws.Range("D94").Value = -12600
$ws.Range("E94").Value = -16600
$ws.Range("F94").Value = -8200

$ws.Range("D100").Value = 77400
$ws.Range("E100").Value = -12100
$ws.Range("F100").Value = 15200
$ws.Range("G100").Value = 40400

$ws.Range("D102").Value = 73000
$ws.Range("E102").Value = -35300
$ws.Range("G102").Value = 32900
